$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Build the new "2022-Q4" sheet by duplicating "2022-Q1" (so it inherits
#    the same header row / column-A styling), inserting it right before
#    "2022-Q1", then overwriting the cell values with the 2022-Q4 figures.
#    This also naturally pushes "2022-Q1", "2021-Q4" and "2020-Q4" one slot
#    later in the tab order, which is exactly what the diff wants.
# ---------------------------------------------------------------------------
$sheetQ1 = $wb.Worksheets.Item("2022-Q1")
$sheetQ1.Copy($sheetQ1)
$newSheet = $wb.Worksheets.Item("2022-Q1 (2)")
$newSheet.Name = "2022-Q4"

# The duplicated sheet has 3 data rows (2,3,4); 2022-Q4 needs 6, so grow it
# by copying the last data row's formatting down for the 3 extra rows.
$newSheet.Rows.Item(4).Copy()
$newSheet.Rows.Item(5).Insert()
$newSheet.Rows.Item(4).Copy()
$newSheet.Rows.Item(5).Insert()
$newSheet.Rows.Item(4).Copy()
$newSheet.Rows.Item(5).Insert()

# Row 2
$newSheet.Cells.Item(2,1).Value = 0
$newSheet.Cells.Item(2,2).Value = "'501030"
$newSheet.Cells.Item(2,3).Value = "汇添富中证环境治理指数（LOF）A"
$newSheet.Cells.Item(2,4).Value = "'2.99"
$newSheet.Cells.Item(2,5).Value = "'92.14"
$newSheet.Cells.Item(2,6).Value = "'1.93"
$newSheet.Cells.Item(2,7).Value = "'0.0577"
$newSheet.Cells.Item(2,8).Value = 7

# Row 3
$newSheet.Cells.Item(3,1).Value = 1
$newSheet.Cells.Item(3,2).Value = "'164908"
$newSheet.Cells.Item(3,3).Value = "交银施罗德中证环境治理指数（LOF）"
$newSheet.Cells.Item(3,4).Value = "'1.55"
$newSheet.Cells.Item(3,5).Value = "'93.92"
$newSheet.Cells.Item(3,6).Value = "'1.98"
$newSheet.Cells.Item(3,7).Value = "'0.0307"
$newSheet.Cells.Item(3,8).Value = 7

# Row 4
$newSheet.Cells.Item(4,1).Value = 2
$newSheet.Cells.Item(4,2).Value = "'501031"
$newSheet.Cells.Item(4,3).Value = "汇添富中证环境治理指数（LOF）C"
$newSheet.Cells.Item(4,4).Value = "'1.38"
$newSheet.Cells.Item(4,5).Value = "'92.14"
$newSheet.Cells.Item(4,6).Value = "'1.93"
$newSheet.Cells.Item(4,7).Value = "'0.0266"
$newSheet.Cells.Item(4,8).Value = 7

# Row 5
$newSheet.Cells.Item(5,1).Value = 3
$newSheet.Cells.Item(5,2).Value = "'000892"
$newSheet.Cells.Item(5,3).Value = "九泰天宝灵活配置混合A"
$newSheet.Cells.Item(5,4).Value = "'0.06"
$newSheet.Cells.Item(5,5).Value = "'94.55"
$newSheet.Cells.Item(5,6).Value = "'4.25"
$newSheet.Cells.Item(5,7).Value = "'0.0026"
$newSheet.Cells.Item(5,8).Value = 8

# Row 6
$newSheet.Cells.Item(6,1).Value = 4
$newSheet.Cells.Item(6,2).Value = "'013413"
$newSheet.Cells.Item(6,3).Value = "交银施罗德中证环境治理指数（LOF）C"
$newSheet.Cells.Item(6,4).Value = "'0.11"
$newSheet.Cells.Item(6,5).Value = "'93.92"
$newSheet.Cells.Item(6,6).Value = "'1.98"
$newSheet.Cells.Item(6,7).Value = "'0.0022"
$newSheet.Cells.Item(6,8).Value = 7

# Row 7
$newSheet.Cells.Item(7,1).Value = 5
$newSheet.Cells.Item(7,2).Value = "'002028"
$newSheet.Cells.Item(7,3).Value = "九泰天宝灵活配置混合C"
$newSheet.Cells.Item(7,4).Value = "'0.00"
$newSheet.Cells.Item(7,5).Value = "'94.55"
$newSheet.Cells.Item(7,6).Value = "'4.25"
$newSheet.Cells.Item(7,7).Value = 0
$newSheet.Cells.Item(7,8).Value = 8

# ---------------------------------------------------------------------------
# 2) Update the "总计" (summary) sheet: insert a new data row for 2022-Q4
#    right after the header, pushing the existing 2022-Q1 / 2021-Q4 / 2020-Q4
#    rows down by one, and renumber the index column.
# ---------------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")

$wsTotal.Rows.Item(2).Copy()
$wsTotal.Rows.Item(3).Insert()

$wsTotal.Cells.Item(2,1).Value = 0
$wsTotal.Cells.Item(2,2).Value = "2022-Q4"
$wsTotal.Cells.Item(2,3).Value = 6
$wsTotal.Cells.Item(2,4).Value = 0.12

$wsTotal.Cells.Item(3,1).Value = 1
$wsTotal.Cells.Item(4,1).Value = 2
$wsTotal.Cells.Item(5,1).Value = 3
